$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ASSOCIATION" / "association" / "SLOT" property column (column E) was
# removed from the dataloader test sheet. Deleting the entire column shifts
# the subsequent columns (F, G) left by one and removes the now-unused
# shared strings automatically.
$ws.Range("E1").EntireColumn.Delete()

# Leave the selection on the (new) column E, matching the saved view state.
$ws.Range("E1").EntireColumn.Select()
